$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell A1: "id" -> "page"
$ws.Range("A1").Value = "page"

# A3: now holds a numeric value
$ws.Range("A3").Value = 9999

# A4: now holds a text value "a", right aligned (matches style index 6)
$ws.Range("A4").Value = "a"
$ws.Range("A4").HorizontalAlignment = -4152
